# "added pay bills section"
# Adds a new "bill-names" worksheet (a lookup list of bill categories) right
# after the existing "report-sheet", and appends five new expense/income
# rows (12-16) to "report-sheet" covering electricity/decoration bills and
# a couple of stock buy/sell transactions.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- new "bill-names" sheet, inserted right after report-sheet ---------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "bill-names"

$ws2.Range("A1").Value = "Bill Names"
$ws2.Range("A2").Value = "Electricity Bill"
$ws2.Range("A3").Value = "Water Bill"
$ws2.Range("A4").Value = "Decoration Fee"
$ws2.Range("A5").Value = "Food cost"
$ws2.Range("A6").Value = "Others"
$ws2.Columns.Item(1).ColumnWidth = 13.75

# --- new rows appended to report-sheet ----------------------------------
$rows = @(
    @{ Row = 12; Date = 44864.74380736111;  Name = "Electricity Bill";                 Amount = -122.0 },
    @{ Row = 13; Date = 44864.76265318287;  Name = "Bought (feemicon pill) X 10000";   Amount = -600000.0 },
    @{ Row = 14; Date = 44864.763558645835; Name = "Bought (Third Party PD - 1) X 3";  Amount = -36.0 },
    @{ Row = 15; Date = 44864.764227037034; Name = "Decoration Fee";                   Amount = -100.0 },
    @{ Row = 16; Date = 44864.76562412037;  Name = "Sold (feemicon pill) X 10";        Amount = 696.9 }
)

foreach ($r in $rows) {
    $dateCell = $ws1.Range("A$($r.Row)")
    $dateCell.Value = $r.Date
    $dateCell.NumberFormat = "mm-dd-yy"

    $ws1.Range("B$($r.Row)").Value = $r.Name
    $ws1.Range("C$($r.Row)").Value = $r.Amount
}

# --- selection state, matching a user who just finished typing row 12 ---
$ws2.Range("A15").Select()
$ws1.Select()
$ws1.Rows.Item(12).Select()
